$d = $word.ActiveDocument
$rsquo = [char]0x2019

# -----------------------------------------------------------------------
# Step 1: Remove the "Meta description" paragraph that currently sits right
# after the title heading (paragraph 2): "Meta description: Discover the
# chemistry-themed slot game with cluster pay and cascade reels. Play Dr
# Reactive's Laboratory free and win big with bonus rounds and progressive
# jackpots."
# -----------------------------------------------------------------------
$metaPara = $d.Paragraphs(2)
$metaPara.Range.Delete()

# -----------------------------------------------------------------------
# Step 2: Insert a new bold paragraph reading "Play Dr Reactive's Laboratory
# Free Slot Game" immediately before the document's final paragraph (the one
# that used to hold the "Create a feature image..." image-generation
# prompt). We do this by inserting a raw OOXML <w:p> fragment at a
# non-boundary offset inside the second-to-last paragraph, which causes the
# new paragraph to be added right after it (i.e. right before the last
# paragraph) without disturbing that paragraph's own text.
# -----------------------------------------------------------------------
$total = $d.Paragraphs.Count
$anchorPara = $d.Paragraphs($total - 1)
$anchorRange = $anchorPara.Range
$insertPoint = $d.Range($anchorRange.Start + 1, $anchorRange.Start + 1)
$null = $insertPoint.InsertXML("<w:p xmlns:w='http://schemas.openxmlformats.org/wordprocessingml/2006/main'><w:r/><w:r><w:rPr><w:b/></w:rPr><w:t>Play Dr Reactive&#8217;s Laboratory Free Slot Game</w:t></w:r></w:p>")

# -----------------------------------------------------------------------
# Step 3: Replace the text of the (now shifted) final paragraph -- which
# still carries its original italic run formatting -- with the new meta
# description copy, dropping the "Meta description: " label and old image
# prompt text entirely.
# -----------------------------------------------------------------------
$total = $d.Paragraphs.Count
$lastPara = $d.Paragraphs($total)
$lastRange = $lastPara.Range
$lastTextRange = $d.Range($lastRange.Start, $lastRange.End - 1)
$lastTextRange.Text = "Discover the chemistry-themed slot game with cluster pay and cascade reels. Play Dr Reactive" + $rsquo + "s Laboratory free and win big with bonus rounds and progressive jackpots."

Write-Host "Edit complete. Total paragraphs:" $d.Paragraphs.Count
